# Update countries & provincias Spain
# Applies the 14-Apr-2020 15:22 -> 15:52 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 15:52"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 587815
$ws.Range("C4").Value = 874
$ws.Range("D4").Value = 37315
$ws.Range("E4").Value = 526846
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 23654

# Reino Unido (row 9)
$ws.Range("B9").Value = 93873
$ws.Range("C9").Value = 5252
$ws.Range("E9").Value = 81422
$ws.Range("G9").Value = 778
$ws.Range("H9").Value = 12107

# Brasil (row 17)
$ws.Range("B17").Value = 23870
$ws.Range("C17").Value = 440
$ws.Range("E17").Value = 19536

# Noruega (row 32)
$ws.Range("B32").Value = 6623
$ws.Range("C32").Value = 20
$ws.Range("E32").Value = 6452
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 139

# Serbia (row 43)
$ws.Range("F43").Value = 131

# Argentina / Sudafrica swap ranking (rows 54-55): Sudafrica moves up to 54
# with refreshed totals, Argentina drops to 55 keeping its prior totals.
$ws.Range("A54").Value = "Sudafrica"
$ws.Range("B54").Value = 2415
$ws.Range("C54").Value = 143
$ws.Range("D54").Value = 410
$ws.Range("E54").Value = 1978
$ws.Range("F54").Value = 7
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 27

$ws.Range("A55").Value = "Argentina"
$ws.Range("B55").Value = 2277
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 559
$ws.Range("E55").Value = 1617
$ws.Range("F55").Value = 83
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 101

# New entry "Republica de Yibuti" enters the table at row 106 with fresh
# data, pushing Bolivia..Senegal down one row each (rows 106-111); the old
# Yibuti row (111) is absorbed by Senegal's shifted data.
$ws.Range("A106").Value = "Republica de Yibuti"
$ws.Range("B106").Value = 363
$ws.Range("C106").Value = 65
$ws.Range("D106").Value = 53
$ws.Range("E106").Value = 308
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 2

$ws.Range("A107").Value = "Bolivia"
$ws.Range("B107").Value = 354
$ws.Range("C107").Value = 24
$ws.Range("D107").Value = 6
$ws.Range("E107").Value = 320
$ws.Range("F107").Value = 3
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 28

$ws.Range("A108").Value = "Nigeria"
$ws.Range("B108").Value = 343
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 91
$ws.Range("E108").Value = 242
$ws.Range("F108").Value = 2
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 10

$ws.Range("A109").Value = "Mauricio"
$ws.Range("B109").Value = 324
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 42
$ws.Range("E109").Value = 273
$ws.Range("F109").Value = 3
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 9

$ws.Range("A110").Value = "Estado de Palestina"
$ws.Range("B110").Value = 308
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 62
$ws.Range("E110").Value = 244
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 2

$ws.Range("A111").Value = "Senegal"
$ws.Range("B111").Value = 299
$ws.Range("C111").Value = 8
$ws.Range("D111").Value = 183
$ws.Range("E111").Value = 114
$ws.Range("F111").Value = 1
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 2

# Consejo Danes para los Refugiados (row 116)
$ws.Range("B116").Value = 241
$ws.Range("C116").Value = 6
$ws.Range("D116").Value = 20
$ws.Range("E116").Value = 201

# Suazilandia (row 185): activos/recuperados swap
$ws.Range("D185").Value = 8
$ws.Range("E185").Value = 7
